$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "release/8.0.7"
$ws.Range("B10").Value = "X"
$ws.Range("C10").Value = "X"
$ws.Range("D10").Value = "X"
$ws.Range("E10").Value = "X"
